# The commit inserts a new data row (a new weekly price observation) right
# before the current row 72 ("Femacal de La Calera" / "Haba"), pushing all
# subsequent rows (72-187) down by one (to 73-188). The new row re-uses the
# same Mercado/Región/Codreg/Categoría/Variedad/Calidad/Unidad/Origen/Kg
# metadata as its neighbours but carries its own Fecha, Volumen,
# Precio máximo, Precio promedio ponderado and Precio $/Kg values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 72; Excel automatically shifts the
# existing rows 72..187 down to 73..188 (and extends the used range to
# A1:R188), carrying along styles/number-formats of the row below.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new observation.
$ws.Cells.Item(72, 1).Value = 3
$ws.Cells.Item(72, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(72, 3).Value = "Coquimbo"
$ws.Cells.Item(72, 4).Value = 44799
$ws.Cells.Item(72, 5).Value = 5
$ws.Cells.Item(72, 6).Value = 100112026
$ws.Cells.Item(72, 7).Value = "Haba"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 55
$ws.Cells.Item(72, 11).Value = 14000
$ws.Cells.Item(72, 12).Value = 14000
$ws.Cells.Item(72, 13).Value = 14000
$ws.Cells.Item(72, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(72, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(72, 16).Value = 560
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = "Hortaliza"
